$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new row for "united states" in its sorted position (before current row 21, "vietnam")
$ws.Rows.Item(21).Insert()
$ws.Range("A21").Value = "united states"

# Insert header row at top
$ws.Rows.Item(1).Insert()
$ws.Range("A1").Value = "Country1"
$ws.Range("B1").Value = "Country2"

for ($r = 2; $r -le 23; $r++) {
    $ws.Cells.Item($r, 2).Formula = "=PROPER(A$r)"
}
